# Generate Report for Handoff
#
# A new handoff batch (XLIFF generation) was run for the files that were
# "Ready for handoff" (and not already fully handed back). This updates:
#   - Overview sheet: "Latest HO Xliff Generate Date" for the affected rows
#   - zh-cn sheet:     "Priority" -> "ht" and "Latest Handoff Datetime" for the affected rows
#   - de-de sheet:     "Priority" -> "ht" and "Latest Handoff Datetime" for the affected rows
#
# The files aa7bb647-... (already handed back) and the two "In Translation"
# files (4a7f418f-..., 6b2cd2e2-...) are not part of this handoff batch.

$wb = $excel.ActiveWorkbook

$rows = @(4, 5, 6, 7, 9, 10)

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-10-25 03:33:53"
}

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-10-25 03:33:40"
}

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-10-25 03:33:53"
}
